# Fruta / hortaliza, semanal
# Insert a new weekly record at row 34 (pushing the existing rows 34-145
# down to 35-146) and populate it with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34:145 down one row, creating a blank row 34.
$ws.Rows("34:34").Insert()

# Populate the new row 34 with this week's data.
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = "2022-05-26"
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112040
$ws.Range("G34").Value = "Cilantro"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 2300
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1750
$ws.Range("N34").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O34").Value = "Provincia del Elquí"
$ws.Range("P34").Value = 1167
$ws.Range("Q34").Value = 1.5
$ws.Range("R34").Value = "Hortaliza"
